$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (runs, balls, fours, sixes) for rows 2-12 after the cyclic
# permutation of row data described by the diff.
$data = @{
    2  = @("67","45","5","4")
    3  = @("29","19","4","1")
    4  = @("20","10","1","1")
    5  = @("41","26","5","0")
    6  = @("9","13","1","0")
    7  = @("57","39","7","0")
    8  = @("22","12","0","2")
    9  = @("50","44","2","2")
    10 = @("8","14","0","0")
    11 = @("11","10","1","0")
    12 = @("3","5","0","0")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").NumberFormat = "@"
    $ws.Range("F$row").Value = $vals[3]
}
